$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 988.24
$ws.Cells.Item(17, 10).Value = 946.1875
$ws.Cells.Item(17, 12).Value = 2838.5625
$ws.Cells.Item(17, 14).Value = -3174.5625

$ws.Cells.Item(32, 8).Value = 1632.4166
$ws.Cells.Item(32, 10).Value = 1222.5
$ws.Cells.Item(32, 12).Value = 1222.5
$ws.Cells.Item(32, 14).Value = -1874.5

$ws.Cells.Item(53, 8).Value = 308.58334
$ws.Cells.Item(53, 9).Value = 175.57143
$ws.Cells.Item(53, 10).Value = 494.8
$ws.Cells.Item(53, 11).Value = 175.57143
$ws.Cells.Item(53, 12).Value = 494.8
$ws.Cells.Item(53, 13).Value = 461.42857
$ws.Cells.Item(53, 14).Value = -1768.8

$ws.Cells.Item(86, 8).Value = 7490.421
$ws.Cells.Item(86, 9).Value = 7544.091
$ws.Cells.Item(86, 10).Value = 7416.625
$ws.Cells.Item(86, 11).Value = 7544.091
$ws.Cells.Item(86, 12).Value = 7416.625
$ws.Cells.Item(86, 13).Value = -6421.091
$ws.Cells.Item(86, 14).Value = -9662.625

$ws.Cells.Item(89, 8).Value = 7490.421
$ws.Cells.Item(89, 9).Value = 7544.091
$ws.Cells.Item(89, 10).Value = 7416.625
$ws.Cells.Item(89, 11).Value = 37720.455
$ws.Cells.Item(89, 12).Value = 37083.125
$ws.Cells.Item(89, 13).Value = -32104.455
$ws.Cells.Item(89, 14).Value = -48315.125

$ws.Cells.Item(112, 8).Value = 2253.7693
$ws.Cells.Item(112, 10).Value = 2209.4546
$ws.Cells.Item(112, 12).Value = 6628.3638
$ws.Cells.Item(112, 14).Value = -8844.363799999999

$ws.Cells.Item(132, 8).Value = 6948.7427
$ws.Cells.Item(132, 9).Value = 7038.2354
$ws.Cells.Item(132, 11).Value = 21114.7062
$ws.Cells.Item(132, 13).Value = -18584.7062

$ws.Cells.Item(138, 8).Value = 7852.0415
$ws.Cells.Item(138, 9).Value = 12675
$ws.Cells.Item(138, 11).Value = 38025
$ws.Cells.Item(138, 13).Value = -32885

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 8000
$ws.Cells.Item(25, 10).Value = 9000
$ws.Cells.Item(25, 12).Value = 9000
$ws.Cells.Item(25, 14).Value = -9804

$ws.Cells.Item(74, 8).Value = 973200.25
$ws.Cells.Item(74, 9).Value = 1118803.2
$ws.Cells.Item(74, 11).Value = 1118803.2
$ws.Cells.Item(74, 13).Value = -1117929.2

$ws.Cells.Item(77, 8).Value = 973200.25
$ws.Cells.Item(77, 9).Value = 1118803.2
$ws.Cells.Item(77, 11).Value = 5594016
$ws.Cells.Item(77, 13).Value = -5589648

$ws.Cells.Item(97, 8).Value = 1041.3077
$ws.Cells.Item(97, 9).Value = 1041.3077
$ws.Cells.Item(97, 11).Value = 1041.3077
$ws.Cells.Item(97, 13).Value = -545.3077000000001

$ws.Cells.Item(110, 8).Value = 1964.4546
$ws.Cells.Item(110, 9).Value = 1899.75
$ws.Cells.Item(110, 10).Value = 2001.4286
$ws.Cells.Item(110, 11).Value = 1899.75
$ws.Cells.Item(110, 12).Value = 2001.4286
$ws.Cells.Item(110, 13).Value = 145.25
$ws.Cells.Item(110, 14).Value = -6091.4286

$ws.Cells.Item(132, 8).Value = 7178.2856
$ws.Cells.Item(132, 9).Value = 3967.3333
$ws.Cells.Item(132, 10).Value = 8054
$ws.Cells.Item(132, 11).Value = 11901.9999
$ws.Cells.Item(132, 12).Value = 24162
$ws.Cells.Item(132, 13).Value = -9371.999899999999
$ws.Cells.Item(132, 14).Value = -29222

$ws.Cells.Item(134, 8).Value = 80000
$ws.Cells.Item(134, 10).Value = 80000
$ws.Cells.Item(134, 12).Value = 80000
$ws.Cells.Item(134, 14).Value = -90140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 1628
$ws.Cells.Item(36, 9).Value = 1628
$ws.Cells.Item(36, 11).Value = 1628
$ws.Cells.Item(36, 13).Value = -1094

$ws.Cells.Item(86, 8).Value = 4249.75
$ws.Cells.Item(86, 9).Value = 5800
$ws.Cells.Item(86, 10).Value = 2699.5
$ws.Cells.Item(86, 11).Value = 5800
$ws.Cells.Item(86, 12).Value = 2699.5
$ws.Cells.Item(86, 13).Value = -4677
$ws.Cells.Item(86, 14).Value = -4945.5

$ws.Cells.Item(89, 8).Value = 4249.75
$ws.Cells.Item(89, 9).Value = 5800
$ws.Cells.Item(89, 10).Value = 2699.5
$ws.Cells.Item(89, 11).Value = 29000
$ws.Cells.Item(89, 12).Value = 13497.5
$ws.Cells.Item(89, 13).Value = -23384
$ws.Cells.Item(89, 14).Value = -24729.5

$ws.Cells.Item(94, 8).Value = 3763.238
$ws.Cells.Item(94, 9).Value = 3106.7896
$ws.Cells.Item(94, 11).Value = 3106.7896
$ws.Cells.Item(94, 13).Value = -2655.7896

$ws.Cells.Item(107, 8).Value = 1507.6086
$ws.Cells.Item(107, 9).Value = 1192
$ws.Cells.Item(107, 10).Value = 1710.5
$ws.Cells.Item(107, 11).Value = 1192
$ws.Cells.Item(107, 12).Value = 1710.5
$ws.Cells.Item(107, 13).Value = 728
$ws.Cells.Item(107, 14).Value = -5550.5

$ws.Cells.Item(135, 8).Value = 90000
$ws.Cells.Item(135, 10).Value = 90000
$ws.Cells.Item(135, 12).Value = 90000
$ws.Cells.Item(135, 14).Value = -100140

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 5786497.5
$ws.Cells.Item(4, 9).Value = 10800313
$ws.Cells.Item(4, 11).Value = 32400939
$ws.Cells.Item(4, 13).Value = -32400827

$ws.Cells.Item(9, 8).Value = 203900.4
$ws.Cells.Item(9, 9).Value = 1000001
$ws.Cells.Item(9, 10).Value = 4875.25
$ws.Cells.Item(9, 11).Value = 3000003
$ws.Cells.Item(9, 12).Value = 14625.75
$ws.Cells.Item(9, 13).Value = -2999779
$ws.Cells.Item(9, 14).Value = -15073.75

$ws.Cells.Item(80, 8).Value = 14333.333
$ws.Cells.Item(80, 10).Value = 18000
$ws.Cells.Item(80, 12).Value = 54000
$ws.Cells.Item(80, 14).Value = -55872

$ws.Cells.Item(83, 8).Value = 14333.333
$ws.Cells.Item(83, 10).Value = 18000
$ws.Cells.Item(83, 12).Value = 162000
$ws.Cells.Item(83, 14).Value = -171360

$ws.Cells.Item(109, 8).Value = 5895.706
$ws.Cells.Item(109, 9).Value = 432.5
$ws.Cells.Item(109, 10).Value = 7576.6924
$ws.Cells.Item(109, 11).Value = 1297.5
$ws.Cells.Item(109, 12).Value = 22730.0772
$ws.Cells.Item(109, 13).Value = -257.5
$ws.Cells.Item(109, 14).Value = -24810.0772

$ws.Cells.Item(131, 8).Value = 6267.9653
$ws.Cells.Item(131, 9).Value = 1343.5454
$ws.Cells.Item(131, 11).Value = 4030.6362
$ws.Cells.Item(131, 13).Value = 1009.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10411.088
$ws.Cells.Item(70, 9).Value = 11636.741
$ws.Cells.Item(70, 10).Value = 5683.5713
$ws.Cells.Item(70, 11).Value = 11636.741
$ws.Cells.Item(70, 12).Value = 5683.5713
$ws.Cells.Item(70, 13).Value = -11366.741
$ws.Cells.Item(70, 14).Value = -6223.5713

$ws.Cells.Item(73, 8).Value = 10411.088
$ws.Cells.Item(73, 9).Value = 11636.741
$ws.Cells.Item(73, 10).Value = 5683.5713
$ws.Cells.Item(73, 11).Value = 11636.741
$ws.Cells.Item(73, 12).Value = 5683.5713
$ws.Cells.Item(73, 13).Value = -10700.741
$ws.Cells.Item(73, 14).Value = -7555.5713

$ws.Cells.Item(132, 8).Value = 44715.855
$ws.Cells.Item(132, 9).Value = 34670.668
$ws.Cells.Item(132, 11).Value = 104012.004
$ws.Cells.Item(132, 13).Value = -101482.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2640
$ws.Cells.Item(7, 9).Value = 2640
$ws.Cells.Item(7, 11).Value = 2640
$ws.Cells.Item(7, 13).Value = -2528

$ws.Cells.Item(46, 8).Value = 3708.25
$ws.Cells.Item(46, 9).Value = 1611.1111
$ws.Cells.Item(46, 11).Value = 1611.1111
$ws.Cells.Item(46, 13).Value = -1423.1111

$ws.Cells.Item(122, 8).Value = 3736
$ws.Cells.Item(122, 10).Value = 4249.2856
$ws.Cells.Item(122, 12).Value = 12747.8568
$ws.Cells.Item(122, 14).Value = -17647.8568

$ws.Cells.Item(126, 8).Value = 2640
$ws.Cells.Item(126, 9).Value = 2640
$ws.Cells.Item(126, 11).Value = 7920
$ws.Cells.Item(126, 13).Value = -5450

$ws.Cells.Item(134, 8).Value = 93436.5
$ws.Cells.Item(134, 10).Value = 93436.5
$ws.Cells.Item(134, 12).Value = 93436.5
$ws.Cells.Item(134, 14).Value = -103576.5

$ws.Cells.Item(136, 8).Value = 11907887
$ws.Cells.Item(136, 9).Value = 6947312.5
$ws.Cells.Item(136, 11).Value = 20841937.5
$ws.Cells.Item(136, 13).Value = -20839387.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3327.4285
$ws.Cells.Item(126, 9).Value = 3358.4
$ws.Cells.Item(126, 10).Value = 3250
$ws.Cells.Item(126, 11).Value = 10075.2
$ws.Cells.Item(126, 12).Value = 9750
$ws.Cells.Item(126, 13).Value = -7605.200000000001
$ws.Cells.Item(126, 14).Value = -14690
